# Add new metamorphic rock entries to the language/localization sheet.
# This mirrors the existing "Name / Desc" row-pair pattern already used
# for igneous and sedimentary rocks further up the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is: (Key, Value, VoiceDuration)
$newRows = @(
    @("amphibolite",          "Amphibolite",        1),
    @("amphiboliteDesc",      "Metamorphic rock.",   5),
    @("anthraciteCoal",       "Anthracite Coal",     2),
    @("anthraciteCoalDesc",   "Metamorphic rock.",   5),
    @("gneiss",               "Gneiss",              0.5),
    @("gneissDesc",           "Metamorphic rock.",   5),
    @("marble",               "Marble",              0.6),
    @("marbleDesc",           "Metamorphic rock.",   5),
    @("metaconglomerate",     "Metaconglomerate",    1.5),
    @("metaconglomerateDesc", "Metamorphic rock.",   5),
    @("phyllite",             "Phyllite",            0.5),
    @("phylliteDesc",         "Metamorphic rock.",   5),
    @("quartzite",            "Quartzite",           1),
    @("quartziteDesc",        "Metamorphic rock.",   5),
    @("schist",               "Schist",              0.5),
    @("schistDesc",           "Metamorphic rock.",   5),
    @("slate",                "Slate",               0.5),
    @("slateDesc",            "Metamorphic rock.",   5)
)

$startRow = 52
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
}

# Update the view so it reflects the author's last scroll/selection position.
$ws.Application.ActiveWindow.ScrollRow = 50
$ws.Range("A70").Select()
